$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-CellText $ws "D2" "60.840.10"
Set-CellText $ws "E2" "  +0.24%  "
Set-CellText $ws "D3" "2.590.55"
Set-CellText $ws "E3" "  -0.09%  "
Set-CellText $ws "E4" "  +0.03%  "
Set-CellText $ws "D5" "522.58"
Set-CellText $ws "E5" "  +2.73%  "
Set-CellText $ws "D6" "154.00"
Set-CellText $ws "E6" "  +0.11%  "
Set-CellText $ws "E7" "  +0.07%  "
Set-CellText $ws "E8" "  +2.48%  "
Set-CellText $ws "E9" "  +1.87%  "
Set-CellText $ws "E10" "  +0.99%  "
Set-CellText $ws "E11" "  -0.48%  "
Set-CellText $ws "E12" "  +1.25%  "
Set-CellText $ws "D13" "3.046.65"
Set-CellText $ws "E13" "  +0.14%  "
Set-CellText $ws "D14" "60.853.90"
Set-CellText $ws "E14" "  +0.40%  "
Set-CellText $ws "E15" "  -0.12%  "
Set-CellText $ws "E16" "  -0.15%  "
Set-CellText $ws "D17" "2.598.57"
Set-CellText $ws "E17" "  +0.28%  "
Set-CellText $ws "E18" "  -0.90%  "
Set-CellText $ws "D19" "352.45"
Set-CellText $ws "E19" "  +1.82%  "
Set-CellText $ws "D20" "10.56"
Set-CellText $ws "E20" "  +1.03%  "
Set-CellText $ws "E21" "  +1.15%  "
Set-CellText $ws "E22" "  +0.16%  "
Set-CellText $ws "D23" "60.77"
Set-CellText $ws "E23" "  +1.20%  "
Set-CellText $ws "E24" "  +1.25%  "
Set-CellText $ws "E25" "  -0.76%  "
Set-CellText $ws "D26" "2.710.40"
Set-CellText $ws "E26" "  +0.31%  "
Set-CellText $ws "E27" "  +0.24%  "
Set-CellText $ws "E28" "  -0.55%  "
Set-CellText $ws "E29" "  -0.28%  "
Set-CellText $ws "E30" "  +0.03%  "
Set-CellText $ws "E31" "  +10.33%  "
Set-CellText $ws "D32" "19.33"
Set-CellText $ws "E32" "  -0.29%  "
Set-CellText $ws "E33" "  +2.42%  "
Set-CellText $ws "D34" "148.02"
Set-CellText $ws "E34" "  -3.78%  "
Set-CellText $ws "E35" "  +3.85%  "
Set-CellText $ws "D36" "0.932"
Set-CellText $ws "E36" "  +8.87%  "
Set-CellText $ws "E37" "  +0.46%  "
Set-CellText $ws "E38" "  +1.47%  "
Set-CellText $ws "E39" "  -0.75%  "
Set-CellText $ws "D40" "3.80"
Set-CellText $ws "E40" "  +0.94%  "
Set-CellText $ws "D41" "36.46"
Set-CellText $ws "E41" "  +1.22%  "
Set-CellText $ws "D42" "288.08"
Set-CellText $ws "E42" "  -2.93%  "
Set-CellText $ws "E43" "  +1.77%  "
Set-CellText $ws "E44" "  +0.40%  "
Set-CellText $ws "D45" "0.0559"
Set-CellText $ws "E45" "  +0.31%  "
Set-CellText $ws "E46" "  +0.13%  "
Set-CellText $ws "E47" "  -1.30%  "
Set-CellText $ws "B48" "RenderToken"
Set-CellText $ws "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D48" "4.88"
Set-CellText $ws "E48" "  +0.22%  "
Set-CellText $ws "B49" "VeChain"
Set-CellText $ws "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D49" "0.0237"
Set-CellText $ws "E49" "  +1.54%  "
Set-CellText $ws "D50" "10.31"
Set-CellText $ws "E50" "  +0.08%  "
Set-CellText $ws "D51" "19.09"
Set-CellText $ws "E51" "  +7.85%  "
